$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data (2020-06-17) continues the daily COVID tracking table one
# more row below the previous last row (96). Copy row 96's formatting down
# into the new row 97 first (date cell keeps the date style, the rest keep
# the centered-number style), then overwrite with the day's real values.
$ws.Range("A96:F96").Copy()
$ws.Range("A97:F97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A97").Value = 43999
$ws.Range("B97").Value = 1294
$ws.Range("C97").Value = 643
$ws.Range("D97").Value = 580
$ws.Range("E97").Value = 303
$ws.Range("F97").Value = 51

# Grow the Excel Table (and its AutoFilter) so the new row is included.
$lo = $ws.ListObjects("Condicion_Pacientes")
$null = $lo.Resize($ws.Range("A1:F97"))

# Match the post-entry view: scrolled down toward the bottom of the table
# with the last cell of the new row selected.
$null = $ws.Range("F97").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 84
$win.ScrollColumn = 4
